# Auto-generated Excel COM-interop script
# Applies market-price / profit value updates across multiple Leve sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) as produced by the scheduled
# market-data runner referenced in the commit message.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 266.2
$ws.Range("I11").Value = 266.2
$ws.Range("K11").Value = 266.2
$ws.Range("M11").Value = -126.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 2746.6
$ws.Range("I69").Value = 1937.6666
$ws.Range("J69").Value = 3960
$ws.Range("K69").Value = 5812.9998
$ws.Range("L69").Value = 11880
$ws.Range("M69").Value = -4938.9998
$ws.Range("N69").Value = -13628

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 2746.6
$ws.Range("I72").Value = 1937.6666
$ws.Range("J72").Value = 3960
$ws.Range("K72").Value = 17438.9994
$ws.Range("L72").Value = 35640
$ws.Range("M72").Value = -13070.9994
$ws.Range("N72").Value = -44376

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2880.6924
$ws.Range("I106").Value = 2955.25
$ws.Range("J106").Value = 1986
$ws.Range("K106").Value = 2955.25
$ws.Range("L106").Value = 1986
$ws.Range("M106").Value = -2324.25
$ws.Range("N106").Value = -3248

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1810
$ws.Range("I137").Value = 1810
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 5430
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -2880
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 484998.66
$ws.Range("I138").Value = 999.8182
$ws.Range("J138").Value = 552390.9
$ws.Range("K138").Value = 2999.4546
$ws.Range("L138").Value = 1657172.7
$ws.Range("M138").Value = 2140.5454
$ws.Range("N138").Value = -1667452.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6151.6313
$ws.Range("I2").Value = 911.1429000000001
$ws.Range("J2").Value = 20825
$ws.Range("K2").Value = 911.1429000000001
$ws.Range("L2").Value = 20825
$ws.Range("M2").Value = -798.1429000000001
$ws.Range("N2").Value = -21051

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2805.7673
$ws.Range("I32").Value = 3084.3333
$ws.Range("K32").Value = 3084.3333
$ws.Range("M32").Value = -2797.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1177.8422
$ws.Range("I74").Value = 960.0769
$ws.Range("J74").Value = 1649.6666
$ws.Range("K74").Value = 960.0769
$ws.Range("L74").Value = 1649.6666
$ws.Range("M74").Value = -86.07690000000002
$ws.Range("N74").Value = -3397.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1177.8422
$ws.Range("I77").Value = 960.0769
$ws.Range("J77").Value = 1649.6666
$ws.Range("K77").Value = 4800.3845
$ws.Range("L77").Value = 8248.333000000001
$ws.Range("M77").Value = -432.3845000000001
$ws.Range("N77").Value = -16984.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H93").Value = 7980
$ws.Range("J93").Value = 7980
$ws.Range("L93").Value = 7980
$ws.Range("N93").Value = -12972

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 6151.6313
$ws.Range("I116").Value = 911.1429000000001
$ws.Range("J116").Value = 20825
$ws.Range("K116").Value = 911.1429000000001
$ws.Range("L116").Value = 20825
$ws.Range("M116").Value = 1382.8571
$ws.Range("N116").Value = -25413

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1685.4286
$ws.Range("I122").Value = 1666.3334
$ws.Range("J122").Value = 1800
$ws.Range("K122").Value = 4999.0002
$ws.Range("L122").Value = 5400
$ws.Range("M122").Value = -2549.0002
$ws.Range("N122").Value = -10300

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2035.3513
$ws.Range("I132").Value = 1736.697
$ws.Range("K132").Value = 5210.090999999999
$ws.Range("M132").Value = -2680.090999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6151.6313
$ws.Range("I3").Value = 911.1429000000001
$ws.Range("J3").Value = 20825
$ws.Range("K3").Value = 911.1429000000001
$ws.Range("L3").Value = 20825
$ws.Range("M3").Value = -797.1429000000001
$ws.Range("N3").Value = -21053

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5244.375
$ws.Range("J86").Value = 4900
$ws.Range("L86").Value = 4900
$ws.Range("N86").Value = -7146

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 5244.375
$ws.Range("J89").Value = 4900
$ws.Range("L89").Value = 24500
$ws.Range("N89").Value = -35732

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1108.909
$ws.Range("I31").Value = 1108.909
$ws.Range("K31").Value = 1108.909
$ws.Range("M31").Value = -813.9090000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1108.909
$ws.Range("I34").Value = 1108.909
$ws.Range("K34").Value = 1108.909
$ws.Range("M34").Value = -906.9090000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H96").Value = 16999
$ws.Range("J96").Value = 16999
$ws.Range("L96").Value = 16999
$ws.Range("N96").Value = -22491

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H97").Value = 24999.5
$ws.Range("J97").Value = 24999.5
$ws.Range("L97").Value = 24999.5
$ws.Range("N97").Value = -26981.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1568.875
$ws.Range("I99").Value = 1507.2858
$ws.Range("K99").Value = 1507.2858
$ws.Range("M99").Value = -9.285800000000108

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1020
$ws.Range("I105").Value = 900
$ws.Range("J105").Value = 1500
$ws.Range("K105").Value = 900
$ws.Range("L105").Value = 1500
$ws.Range("M105").Value = 847
$ws.Range("N105").Value = -4994

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1568.875
$ws.Range("I126").Value = 1507.2858
$ws.Range("K126").Value = 4521.857400000001
$ws.Range("M126").Value = -2051.857400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H127").Value = 39709
$ws.Range("I127").Value = 39709
$ws.Range("K127").Value = 39709
$ws.Range("M127").Value = -34749

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1015277.8
$ws.Range("I4").Value = 449544.5
$ws.Range("J4").Value = 1109566.6
$ws.Range("K4").Value = 1348633.5
$ws.Range("L4").Value = 3328699.8
$ws.Range("M4").Value = -1348521.5
$ws.Range("N4").Value = -3328923.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 413.14285
$ws.Range("I14").Value = 413.14285
$ws.Range("K14").Value = 1239.42855
$ws.Range("M14").Value = -1066.42855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 751.2
$ws.Range("I132").Value = 751.2
$ws.Range("K132").Value = 6760.8
$ws.Range("M132").Value = -4230.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 2000
$ws.Range("I133").Value = 2000
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 6000
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -940
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2159.2285
$ws.Range("I132").Value = 1819.8334
$ws.Range("J132").Value = 4195.6
$ws.Range("K132").Value = 5459.5002
$ws.Range("L132").Value = 12586.8
$ws.Range("M132").Value = -2929.5002
$ws.Range("N132").Value = -17646.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2435
$ws.Range("I7").Value = 2150
$ws.Range("J7").Value = 3005
$ws.Range("K7").Value = 2150
$ws.Range("L7").Value = 3005
$ws.Range("M7").Value = -2038
$ws.Range("N7").Value = -3229

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 35717984
$ws.Range("J122").Value = 4999.5
$ws.Range("L122").Value = 14998.5
$ws.Range("N122").Value = -19898.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2435
$ws.Range("I126").Value = 2150
$ws.Range("J126").Value = 3005
$ws.Range("K126").Value = 6450
$ws.Range("L126").Value = 9015
$ws.Range("M126").Value = -3980
$ws.Range("N126").Value = -13955

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5066.778
$ws.Range("J81").Value = 5606.3125
$ws.Range("L81").Value = 11212.625
$ws.Range("N81").Value = -13334.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 5066.778
$ws.Range("J84").Value = 5606.3125
$ws.Range("L84").Value = 56063.125
$ws.Range("N84").Value = -66671.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 125001580
$ws.Range("I126").Value = 142858700
$ws.Range("K126").Value = 428576100
$ws.Range("M126").Value = -428573630

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 35531.668
$ws.Range("J133").Value = 35531.668
$ws.Range("L133").Value = 35531.668
$ws.Range("N133").Value = -45651.668
